$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '27.406.07'
$ws.Range('E2').Value = '  +1.14%  '
$ws.Range('D3').Value = '1.859.14'
$ws.Range('E3').Value = '  +1.54%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.47%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('E6').Value = '  -0.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4627'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3721'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07318'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8888'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.84%  '
$ws.Range('E11').Value = '  +0.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07819'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.02%  '
$ws.Range('D13').Value = '1.805.57'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.93'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.004'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008975'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('D21').Value = '27.416.53'
$ws.Range('E21').Value = '  +1.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.130'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('D24').Value = '2.086.08'
$ws.Range('E24').Value = '  +3.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.935'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.22'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.065'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '116.43'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08852'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.134'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.39%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7698'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.175'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.516'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.712'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +10.38%  '
$ws.Range('E37').Value = '  +0.65%  '
$ws.Range('E38').Value = '  +0.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05240'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.957'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.086'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5148'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.85%  '
$ws.Range('E43').Value = '  +0.43%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.415'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4815'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.32'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.90%  '
$ws.Range('E47').Value = '  -0.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '103.45'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.97%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.653'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06220'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '65.60'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.94%  '
